# A few more test cases were written/relocated for the "Step3" folder
# (tracked in the external workbook Step3/_Test_Suite_Statistics.xlsx).
# Its Total Test Cases (G4) went 98 -> 108 and its Automated Test Cases
# (G5) went 75 -> 85. Those two external values feed this sheet's
# "Step3" row (row 5) through F5 (=[4]Sheet1!$G$4) and E5
# (=[4]Sheet1!$G$5), so push the refreshed totals through those cells;
# the roll-up formulas in H5/H6/H7 (SUM($F:$F), SUM($E:$E), H6/H5) then
# recompute on their own from the new row-5 numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F5").Formula = "=108"
$ws.Range("E5").Formula = "=85"
